$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (Cadastrado) with new totals
$ws.Range("B2").Value = 115677.51
$ws.Range("B3").Value = 953108.38
$ws.Range("B4").Value = 1786705.19
$ws.Range("B5").Value = 2885974.02
$ws.Range("B6").Value = 4517432.77
$ws.Range("B7").Value = 1774017.75

# Remove the "Sem Cadastro" column entirely
$ws.Range("C1:C7").EntireColumn.Delete()
